$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6 through 23 first (from bottom up to keep row indices stable isn't required since we delete a contiguous block once)
$ws.Range("A6:A23").EntireRow.Delete()

$ws.Range("A2").Value = "('Exalted Angel', ['{4}{W}{W}', 'Creature — Angel', 'Flying', 'Whenever Exalted Angel deals damage, you gain that much life.', 'Morph {2}{W}{W} (You may cast this card face down as a 2/2 creature for {3}. Turn it face up any time for its morph cost.)', '4/5'])"
$ws.Range("A3").Value = "('Grim Lavamancer', ['{R}', 'Creature — Human Wizard', '{R}, {T}, Exile two cards from your graveyard: Grim Lavamancer deals 2 damage to any target.', '1/1'])"
$ws.Range("A4").Value = "('Meddling Mage', ['{W}{U}', 'Creature — Human Wizard', 'As Meddling Mage enters the battlefield, choose a nonland card name.', 'Spells with the chosen name can’t be cast.', '2/2'])"
$ws.Range("A5").Value = "('Pernicious Deed', ['{1}{B}{G}', 'Enchantment', '{X}, Sacrifice Pernicious Deed: Destroy each artifact, creature, and enchantment with converted mana cost X or less.'])"
